# Auto-generated script applying 2023-02-27 crime data update
# across the "Citywide Totals", "By Neighborhood" and per-neighborhood sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 10).Value = 983
$ws.Cells.Item(3, 10).Value = 1065
$ws.Cells.Item(4, 10).Value = 78
$ws.Cells.Item(5, 9).Value = 7591
$ws.Cells.Item(5, 10).Value = 1184
$ws.Cells.Item(6, 6).Value = 1879
$ws.Cells.Item(6, 10).Value = 241
$ws.Cells.Item(7, 10).Value = 78
$ws.Cells.Item(8, 10).Value = 4869
$ws.Cells.Item(9, 10).Value = 1460
$ws.Cells.Item(10, 10).Value = 7683
$ws.Cells.Item(11, 6).Value = 105545
$ws.Cells.Item(11, 9).Value = 110282
$ws.Cells.Item(11, 10).Value = 17641

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 10).Value = 161
$ws.Cells.Item(4, 10).Value = 88
$ws.Cells.Item(6, 10).Value = 175
$ws.Cells.Item(7, 10).Value = 455
$ws.Cells.Item(8, 10).Value = 720
$ws.Cells.Item(9, 10).Value = 93
$ws.Cells.Item(10, 10).Value = 150
$ws.Cells.Item(11, 10).Value = 286
$ws.Cells.Item(13, 10).Value = 28
$ws.Cells.Item(14, 10).Value = 104
$ws.Cells.Item(15, 10).Value = 154
$ws.Cells.Item(16, 10).Value = 119
$ws.Cells.Item(18, 10).Value = 137
$ws.Cells.Item(19, 10).Value = 453
$ws.Cells.Item(20, 10).Value = 288
$ws.Cells.Item(23, 10).Value = 248
$ws.Cells.Item(24, 10).Value = 94
$ws.Cells.Item(26, 10).Value = 38
$ws.Cells.Item(27, 10).Value = 220
$ws.Cells.Item(29, 10).Value = 618
$ws.Cells.Item(30, 6).Value = 308
$ws.Cells.Item(31, 10).Value = 141
$ws.Cells.Item(33, 10).Value = 439
$ws.Cells.Item(34, 10).Value = 152
$ws.Cells.Item(36, 10).Value = 234
$ws.Cells.Item(37, 10).Value = 467
$ws.Cells.Item(38, 10).Value = 15
$ws.Cells.Item(39, 10).Value = 18
$ws.Cells.Item(40, 10).Value = 47
$ws.Cells.Item(41, 10).Value = 85
$ws.Cells.Item(42, 10).Value = 469
$ws.Cells.Item(43, 10).Value = 191
$ws.Cells.Item(44, 10).Value = 198
$ws.Cells.Item(46, 10).Value = 60
$ws.Cells.Item(48, 10).Value = 362
$ws.Cells.Item(49, 10).Value = 214
$ws.Cells.Item(50, 10).Value = 180
$ws.Cells.Item(51, 10).Value = 233
$ws.Cells.Item(52, 10).Value = 265
$ws.Cells.Item(53, 10).Value = 272
$ws.Cells.Item(54, 10).Value = 520
$ws.Cells.Item(55, 10).Value = 175
$ws.Cells.Item(56, 10).Value = 87
$ws.Cells.Item(59, 10).Value = 46
$ws.Cells.Item(60, 9).Value = 768
$ws.Cells.Item(60, 10).Value = 127
$ws.Cells.Item(62, 10).Value = 5
$ws.Cells.Item(63, 10).Value = 224
$ws.Cells.Item(64, 10).Value = 148
$ws.Cells.Item(65, 10).Value = 261
$ws.Cells.Item(66, 10).Value = 108
$ws.Cells.Item(67, 10).Value = 369
$ws.Cells.Item(70, 10).Value = 120
$ws.Cells.Item(71, 10).Value = 70
$ws.Cells.Item(72, 10).Value = 81
$ws.Cells.Item(73, 10).Value = 181
$ws.Cells.Item(75, 10).Value = 70
$ws.Cells.Item(76, 10).Value = 515
$ws.Cells.Item(78, 10).Value = 273
$ws.Cells.Item(79, 10).Value = 411
$ws.Cells.Item(80, 10).Value = 49
$ws.Cells.Item(82, 10).Value = 41
$ws.Cells.Item(83, 10).Value = 312
$ws.Cells.Item(85, 10).Value = 697
$ws.Cells.Item(86, 10).Value = 119
$ws.Cells.Item(87, 10).Value = 62
$ws.Cells.Item(88, 10).Value = 157
$ws.Cells.Item(89, 10).Value = 263
$ws.Cells.Item(90, 10).Value = 263
$ws.Cells.Item(91, 10).Value = 178
$ws.Cells.Item(93, 10).Value = 130
$ws.Cells.Item(94, 10).Value = 373
$ws.Cells.Item(95, 10).Value = 252
$ws.Cells.Item(96, 10).Value = 256
$ws.Cells.Item(97, 10).Value = 259
$ws.Cells.Item(98, 10).Value = 165
$ws.Cells.Item(99, 10).Value = 257
$ws.Cells.Item(100, 10).Value = 40
$ws.Cells.Item(101, 6).Value = 105545
$ws.Cells.Item(101, 9).Value = 110282
$ws.Cells.Item(101, 10).Value = 17641

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(2, 10).Value = 14
$ws.Cells.Item(10, 10).Value = 143
$ws.Cells.Item(11, 10).Value = 263

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 10).Value = 40
$ws.Cells.Item(3, 10).Value = 52
$ws.Cells.Item(5, 10).Value = 65
$ws.Cells.Item(7, 10).Value = 6
$ws.Cells.Item(8, 10).Value = 271
$ws.Cells.Item(9, 10).Value = 46
$ws.Cells.Item(10, 10).Value = 203
$ws.Cells.Item(11, 10).Value = 697

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 10).Value = 21
$ws.Cells.Item(4, 10).Value = 6
$ws.Cells.Item(8, 10).Value = 45
$ws.Cells.Item(10, 10).Value = 115
$ws.Cells.Item(11, 10).Value = 265

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(4, 10).Value = 4
$ws.Cells.Item(10, 10).Value = 135
$ws.Cells.Item(11, 10).Value = 286

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 10).Value = 77
$ws.Cells.Item(5, 10).Value = 32
$ws.Cells.Item(6, 10).Value = 15
$ws.Cells.Item(8, 10).Value = 221
$ws.Cells.Item(9, 10).Value = 71
$ws.Cells.Item(10, 10).Value = 213
$ws.Cells.Item(11, 10).Value = 720

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(10, 10).Value = 122
$ws.Cells.Item(11, 10).Value = 272

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 10).Value = 32
$ws.Cells.Item(8, 10).Value = 165
$ws.Cells.Item(10, 10).Value = 135
$ws.Cells.Item(11, 10).Value = 455

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(10, 10).Value = 115
$ws.Cells.Item(11, 10).Value = 256

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Cells.Item(10, 10).Value = 91
$ws.Cells.Item(11, 10).Value = 120

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Cells.Item(10, 10).Value = 46
$ws.Cells.Item(11, 10).Value = 104

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Cells.Item(6, 6).Value = 9
$ws.Cells.Item(11, 6).Value = 308

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(3, 10).Value = 43
$ws.Cells.Item(8, 10).Value = 174
$ws.Cells.Item(9, 10).Value = 44
$ws.Cells.Item(10, 10).Value = 131
$ws.Cells.Item(11, 10).Value = 467

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(6, 10).Value = 2
$ws.Cells.Item(11, 10).Value = 257

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(6, 10).Value = 9
$ws.Cells.Item(8, 10).Value = 98
$ws.Cells.Item(11, 10).Value = 369

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(5, 10).Value = 15
$ws.Cells.Item(10, 10).Value = 60
$ws.Cells.Item(11, 10).Value = 141

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(6, 10).Value = 6
$ws.Cells.Item(8, 10).Value = 54
$ws.Cells.Item(9, 10).Value = 44
$ws.Cells.Item(11, 10).Value = 261

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(6, 10).Value = 4
$ws.Cells.Item(8, 10).Value = 38
$ws.Cells.Item(11, 10).Value = 220

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(2, 10).Value = 25
$ws.Cells.Item(10, 10).Value = 95
$ws.Cells.Item(11, 10).Value = 312

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(5, 10).Value = 13
$ws.Cells.Item(10, 10).Value = 143
$ws.Cells.Item(11, 10).Value = 439

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(2, 10).Value = 25
$ws.Cells.Item(8, 10).Value = 96
$ws.Cells.Item(10, 10).Value = 77
$ws.Cells.Item(11, 10).Value = 252

$ws = $wb.Worksheets.Item('Pullman')
$ws.Cells.Item(10, 10).Value = 27
$ws.Cells.Item(11, 10).Value = 70

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(2, 10).Value = 33
$ws.Cells.Item(5, 10).Value = 23
$ws.Cells.Item(8, 10).Value = 161
$ws.Cells.Item(10, 10).Value = 108
$ws.Cells.Item(11, 10).Value = 411

$ws = $wb.Worksheets.Item('Oakland')
$ws.Cells.Item(3, 10).Value = 3
$ws.Cells.Item(8, 10).Value = 33
$ws.Cells.Item(11, 10).Value = 70

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Cells.Item(10, 10).Value = 23
$ws.Cells.Item(11, 10).Value = 60

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Cells.Item(2, 10).Value = 4
$ws.Cells.Item(3, 10).Value = 6
$ws.Cells.Item(10, 9).Value = 379
$ws.Cells.Item(10, 10).Value = 51
$ws.Cells.Item(11, 9).Value = 768
$ws.Cells.Item(11, 10).Value = 127

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Cells.Item(2, 10).Value = 4
$ws.Cells.Item(3, 10).Value = 4
$ws.Cells.Item(10, 10).Value = 71
$ws.Cells.Item(11, 10).Value = 148

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Cells.Item(8, 10).Value = 15
$ws.Cells.Item(11, 10).Value = 47

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(10, 10).Value = 306
$ws.Cells.Item(11, 10).Value = 373

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(5, 10).Value = 20
$ws.Cells.Item(8, 10).Value = 41
$ws.Cells.Item(10, 10).Value = 394
$ws.Cells.Item(11, 10).Value = 515

$ws = $wb.Worksheets.Item('North Center')
$ws.Cells.Item(10, 10).Value = 62
$ws.Cells.Item(11, 10).Value = 108

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Cells.Item(10, 10).Value = 30
$ws.Cells.Item(11, 10).Value = 62

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Cells.Item(9, 10).Value = 23
$ws.Cells.Item(10, 10).Value = 40

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Cells.Item(10, 10).Value = 76
$ws.Cells.Item(11, 10).Value = 119

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Cells.Item(8, 10).Value = 26
$ws.Cells.Item(11, 10).Value = 214

$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(9, 10).Value = 15
$ws.Cells.Item(11, 10).Value = 259

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(3, 10).Value = 8
$ws.Cells.Item(10, 10).Value = 79
$ws.Cells.Item(11, 10).Value = 175

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(8, 10).Value = 62
$ws.Cells.Item(9, 10).Value = 39
$ws.Cells.Item(10, 10).Value = 376
$ws.Cells.Item(11, 10).Value = 520

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(10, 10).Value = 94
$ws.Cells.Item(11, 10).Value = 181

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(3, 10).Value = 73
$ws.Cells.Item(6, 10).Value = 9
$ws.Cells.Item(8, 10).Value = 175
$ws.Cells.Item(9, 10).Value = 55
$ws.Cells.Item(10, 10).Value = 191
$ws.Cells.Item(11, 10).Value = 618

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(2, 10).Value = 27
$ws.Cells.Item(5, 10).Value = 33
$ws.Cells.Item(7, 10).Value = 9
$ws.Cells.Item(8, 10).Value = 142
$ws.Cells.Item(9, 10).Value = 52
$ws.Cells.Item(10, 10).Value = 146
$ws.Cells.Item(11, 10).Value = 453

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(8, 10).Value = 48
$ws.Cells.Item(10, 10).Value = 105
$ws.Cells.Item(11, 10).Value = 198

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 10).Value = 31
$ws.Cells.Item(5, 10).Value = 24
$ws.Cells.Item(8, 10).Value = 136
$ws.Cells.Item(11, 10).Value = 469

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(8, 10).Value = 63
$ws.Cells.Item(10, 10).Value = 246
$ws.Cells.Item(11, 10).Value = 362

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Cells.Item(8, 10).Value = 64
$ws.Cells.Item(11, 10).Value = 175

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Cells.Item(3, 10).Value = 5
$ws.Cells.Item(8, 10).Value = 36
$ws.Cells.Item(11, 10).Value = 85

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(8, 10).Value = 93
$ws.Cells.Item(9, 10).Value = 27
$ws.Cells.Item(11, 10).Value = 234

$ws = $wb.Worksheets.Item('Boystown')
$ws.Cells.Item(5, 10).Value = 2
$ws.Cells.Item(10, 10).Value = 28

$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(8, 10).Value = 36
$ws.Cells.Item(10, 10).Value = 78
$ws.Cells.Item(11, 10).Value = 150

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Cells.Item(10, 10).Value = 82
$ws.Cells.Item(11, 10).Value = 119

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(8, 10).Value = 67
$ws.Cells.Item(9, 10).Value = 15
$ws.Cells.Item(10, 10).Value = 142
$ws.Cells.Item(11, 10).Value = 273

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(10, 10).Value = 58
$ws.Cells.Item(11, 10).Value = 154

$ws = $wb.Worksheets.Item('Dunning')
$ws.Cells.Item(10, 10).Value = 37
$ws.Cells.Item(11, 10).Value = 94

$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(10, 10).Value = 92
$ws.Cells.Item(11, 10).Value = 248

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(3, 10).Value = 30
$ws.Cells.Item(8, 10).Value = 74
$ws.Cells.Item(10, 10).Value = 116
$ws.Cells.Item(11, 10).Value = 288

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(5, 10).Value = 10
$ws.Cells.Item(11, 10).Value = 233

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(3, 10).Value = 24
$ws.Cells.Item(10, 10).Value = 43
$ws.Cells.Item(11, 10).Value = 178

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(3, 10).Value = 12
$ws.Cells.Item(8, 10).Value = 122
$ws.Cells.Item(9, 10).Value = 21
$ws.Cells.Item(11, 10).Value = 263

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Cells.Item(10, 10).Value = 79
$ws.Cells.Item(11, 10).Value = 180

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Cells.Item(3, 10).Value = 8
$ws.Cells.Item(8, 10).Value = 41
$ws.Cells.Item(11, 10).Value = 130

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(8, 10).Value = 37
$ws.Cells.Item(10, 10).Value = 33
$ws.Cells.Item(11, 10).Value = 137

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Cells.Item(10, 10).Value = 80
$ws.Cells.Item(11, 10).Value = 87

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(9, 10).Value = 8
$ws.Cells.Item(11, 10).Value = 161

$ws = $wb.Worksheets.Item('Old Town')
$ws.Cells.Item(8, 10).Value = 16
$ws.Cells.Item(11, 10).Value = 81

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(5, 10).Value = 13
$ws.Cells.Item(9, 10).Value = 29
$ws.Cells.Item(11, 10).Value = 191

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Cells.Item(10, 10).Value = 36
$ws.Cells.Item(11, 10).Value = 88

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Cells.Item(10, 10).Value = 34
$ws.Cells.Item(11, 10).Value = 49

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Cells.Item(2, 10).Value = 6
$ws.Cells.Item(11, 10).Value = 152

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Cells.Item(5, 10).Value = 19
$ws.Cells.Item(8, 10).Value = 31
$ws.Cells.Item(10, 10).Value = 86
$ws.Cells.Item(11, 10).Value = 165

$ws = $wb.Worksheets.Item('East Village')
$ws.Cells.Item(10, 10).Value = 14
$ws.Cells.Item(11, 10).Value = 38

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Cells.Item(5, 10).Value = 6
$ws.Cells.Item(10, 10).Value = 41

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Cells.Item(5, 10).Value = 5
$ws.Cells.Item(11, 10).Value = 93

$ws = $wb.Worksheets.Item('Montclare')
$ws.Cells.Item(2, 10).Value = 5
$ws.Cells.Item(11, 10).Value = 46

$ws = $wb.Worksheets.Item('Greektown')
$ws.Cells.Item(8, 10).Value = 12
$ws.Cells.Item(9, 10).Value = 18

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(8, 10).Value = 57
$ws.Cells.Item(10, 10).Value = 61
$ws.Cells.Item(11, 10).Value = 157

$ws = $wb.Worksheets.Item('Grant Park')
$ws.Cells.Item(9, 10).Value = 12
$ws.Cells.Item(10, 10).Value = 15
$ws.Cells.Item(8, 10).ClearContents()

$ws = $wb.Worksheets.Item('Museum Campus')
$ws.Cells.Item(8, 10).Value = 5
$ws.Cells.Item(9, 10).Value = 5
